# Update attendance/interest numbers ("想去人数") in "展览" and "全部类型" sheets
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheets index 1)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 854
$wsExhibit.Range("F6").Value = 2102
$wsExhibit.Range("F7").Value = 189

# Sheet "全部类型" (Worksheets index 4)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 854
$wsAll.Range("F8").Value = 2102
$wsAll.Range("F10").Value = 189
